# Update the three metabolite values in row 20 (Sheet1) to reflect the
# recalculated data used for the new combined-metabolite graphs.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A20").Value = 2576.2845482917992
$ws.Range("B20").Value = 1981.3402408926472
$ws.Range("C20").Value = 1870.1560410687173
